$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the "Personal communication" source lines down by one row ---
# Before: A52=Meri Daushvili.., A53=(blank), A54=Tinatin Ksovreli.., A55=(blank), A56=National Statistics.., A57=(blank)
# After:  A52=(blank), A53=Meri Daushvili.., A54=(blank), A55=Tinatin Ksovreli.., A56=(blank), A57=National Statistics..

$meri = $ws.Range("A52").Value2
$tinatin = $ws.Range("A54").Value2
$natstat = $ws.Range("A56").Value2

$ws.Range("A52").Value = ""
$ws.Range("A53").Value = $meri
$ws.Range("A54").Value = ""
$ws.Range("A55").Value = $tinatin
$ws.Range("A56").Value = ""
$ws.Range("A57").Value = $natstat

# --- Replace the GeoStat source citation text in A61 ---
$ws.Range("A61").Value = '"SMALL AND MEDIUM BUSINESS IN GEORGIA" by Department of Statistics, p. 6. Available at http://geostat.ge/cms/site_images/_files/english/statistika%20eng%202009.pdf'
